$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2937
$ws1.Range("F10").Value = 6876
$ws1.Range("F12").Value = 63
$ws1.Range("F13").Value = 348
$ws1.Range("F15").Value = 1487
$ws1.Range("F19").Value = 650
$ws1.Range("G19").Value = "不可售"
$ws1.Range("F20").Value = 106
$ws1.Range("F23").Value = 175
$ws1.Range("F24").Value = 336
$ws1.Range("F25").Value = 1700
$ws1.Range("F26").Value = 1686
$ws1.Range("F27").Value = 535
$ws1.Range("F31").Value = 1211
$ws1.Range("F32").Value = 137
$ws1.Range("F34").Value = 29
$ws1.Range("F36").Value = 422
$ws1.Range("F37").Value = 10
$ws1.Range("F38").Value = 2456
$ws1.Range("F39").Value = 2713
$ws1.Range("F40").Value = 71
$ws1.Range("F45").Value = 316
$ws1.Range("F47").Value = 166

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 162
$ws2.Range("F13").Value = 3
$ws2.Range("F19").Value = 45
$ws2.Range("F23").Value = 467
$ws2.Range("F24").Value = 42

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 538
$ws3.Range("F8").Value = 2719
$ws3.Range("F9").Value = 1004
$ws3.Range("F10").Value = 913
$ws3.Range("F11").Value = 36
$ws3.Range("F12").Value = 255
$ws3.Range("F13").Value = 1432
$ws3.Range("F14").Value = 7314

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 538
$ws4.Range("F5").Value = 2937
$ws4.Range("F8").Value = 2719
$ws4.Range("F9").Value = 6876
$ws4.Range("F10").Value = 1004
$ws4.Range("F12").Value = 348
$ws4.Range("F13").Value = 162
$ws4.Range("F14").Value = 255
$ws4.Range("C18").Value = "上海·（国际）微缩艺术模玩展-GMHS 2024"
$ws4.Range("D18").Value = "国展路1099号 上海世博展览馆"
$ws4.Range("E18").Value = "2024.08.17 09:00-08.18 17:00"
$ws4.Range("F18").Value = 106
$ws4.Range("G18").Value = 49
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=90343"
$ws4.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202408/L7dY65lZ1722843040285.jpeg"
$ws4.Range("B19").Value = "2024-08-18"
$ws4.Range("C19").Value = "上海·东方PartyNight"
$ws4.Range("D19").Value = "重庆南路308号3楼（近建国中路） 上海MaoLivehouse"
$ws4.Range("E19").Value = "2024.08.18 22:00-08.19 01:00"
$ws4.Range("F19").Value = 180
$ws4.Range("G19").Value = 149
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=89209"
$ws4.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202407/4lZtvl551720680564562.jpeg"
$ws4.Range("C20").Value = "上海·第十四届ACBC动漫游戏盛典（免费活动）"
$ws4.Range("D20").Value = "御北路515号 红星美凯龙浦东沪南商场"
$ws4.Range("E20").Value = "2024.08.18 11:00-08.18 18:00"
$ws4.Range("F20").Value = 1107
$ws4.Range("G20").Value = 80
$ws4.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=90508"
$ws4.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202408/qUObKFor1723124901300.jpeg"
$ws4.Range("B21").Value = "2024-08-24"
$ws4.Range("C21").Value = "上海·HAG 1st live in Shanghai《不眨眼》2024演唱会"
$ws4.Range("D21").Value = "中兴路1683号金融街购物中心三楼L3-27 蜚声LIVE House"
$ws4.Range("E21").Value = "2024.08.24 19:30-08.24 21:30"
$ws4.Range("F21").Value = 58
$ws4.Range("G21").Value = 480
$ws4.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=89977"
$ws4.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202407/iXZNZNM01722243246403.png"
$ws4.Range("C22").Value = "上海·coser动漫展"
$ws4.Range("D22").Value = "海潮路133号B1 JUMP工坊"
$ws4.Range("E22").Value = "2024.08.24 10:00-08.25 17:00"
$ws4.Range("F22").Value = 1700
$ws4.Range("G22").Value = 60
$ws4.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=87347"
$ws4.Range("I22").Value = "//i0.hdslb.com/bfs/openplatform/202406/i6vAgX8I1719311206769.jpeg"
$ws4.Range("C23").Value = "上海·幻梦游戏律动——怪物猎人X最终幻想X塞尔达 燃炸游戏室内乐音乐会"
$ws4.Range("D23").Value = "南京西路1376号 上海商城剧院"
$ws4.Range("E23").Value = "2024.08.24 19:30-08.24 21:00"
$ws4.Range("F23").Value = 163
$ws4.Range("G23").Value = 128
$ws4.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=85461"
$ws4.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202405/7fJJ5GxW1715327101441.jpeg"
$ws4.Range("C24").Value = "上海·火舞之时二次元轰趴"
$ws4.Range("D24").Value = "政通路189号五角场万达广场C栋 元气森林livehouse"
$ws4.Range("E24").Value = "2024.08.24 13:00-08.24 19:00"
$ws4.Range("F24").Value = 35
$ws4.Range("G24").Value = 169
$ws4.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=90450"
$ws4.Range("I24").Value = "//i0.hdslb.com/bfs/openplatform/202408/IFeEcHjF1723087744637.png"
$ws4.Range("C25").Value = "上海·第68届燃梦星辰国潮动漫嘉年华-次元盛典我们在燃梦相遇吧！（免费展）"
$ws4.Range("D25").Value = "云锦路500号(近11号线地铁站5号口) 绿地滨江CLUB"
$ws4.Range("E25").Value = "2024.08.24 13:30-08.24 17:00"
$ws4.Range("F25").Value = 1661
$ws4.Range("G25").Value = 58.8
$ws4.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=89301"
$ws4.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202407/TCo7vHap1721008552929.jpeg"
$ws4.Range("C26").Value = "上海·第五人格同人only 同人展4.0"
$ws4.Range("D26").Value = "吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙"
$ws4.Range("E26").Value = "2024.08.24 10:00-08.25 17:00"
$ws4.Range("F26").Value = 1211
$ws4.Range("G26").Value = 68
$ws4.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=90331"
$ws4.Range("I26").Value = "//i0.hdslb.com/bfs/openplatform/202408/WBPPS7uZ1722582105314.jpeg"
$ws4.Range("C27").Value = "上海·第十六届Redamancy动漫游戏嘉年华"
$ws4.Range("D27").Value = "中山北路3300号 上海JOYPOLIS世嘉都市乐园"
$ws4.Range("F27").Value = 137
$ws4.Range("G27").Value = 60
$ws4.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=87713"
$ws4.Range("I27").Value = "//i0.hdslb.com/bfs/openplatform/202406/NVGDyUdo1718294083363.png"
$ws4.Range("B28").Value = "2024-08-25"
$ws4.Range("C28").Value = "上海·「浪漫主义之夜：肖邦遇见贝多芬」罗赛·罗森博伊姆钢琴独奏音乐会"
$ws4.Range("D28").Value = "丁香路425号(上海科技馆地铁站1号口步行460米) 上海东方艺术中心音乐厅"
$ws4.Range("E28").Value = "2024.08.25 19:30-08.25 21:20"
$ws4.Range("F28").Value = 2
$ws4.Range("G28").Value = 268
$ws4.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=90491"
$ws4.Range("I28").Value = "//i2.hdslb.com/bfs/openplatform/202408/tvC4Hl8h1722858714406.jpeg"
$ws4.Range("F30").Value = 29
$ws4.Range("F34").Value = 467
$ws4.Range("F35").Value = 422
$ws4.Range("F37").Value = 10
$ws4.Range("F38").Value = 2456
$ws4.Range("F39").Value = 2713
$ws4.Range("F40").Value = 71
$ws4.Range("F44").Value = 316
$ws4.Range("F46").Value = 166
